$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.974.72'
$ws.Range('E2').Value = '  +1.66%  '
$ws.Range('D3').Value = '2.998.48'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''561.28'
$ws.Range('E5').Value = '  +1.48%  '
$ws.Range('D6').Value = '''138.03'
$ws.Range('E6').Value = '  +11.41%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +4.58%  '
$ws.Range('D9').Value = '2.987.38'
$ws.Range('E9').Value = '  +2.46%  '
$ws.Range('E10').Value = '  +4.50%  '
$ws.Range('D11').Value = '''4.87'
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('D12').Value = '''0.456'
$ws.Range('E12').Value = '  +3.78%  '
$ws.Range('D13').Value = '''0.0000230'
$ws.Range('E13').Value = '  +7.48%  '
$ws.Range('D14').Value = '''33.77'
$ws.Range('E14').Value = '  +4.93%  '
$ws.Range('E15').Value = '  +2.81%  '
$ws.Range('D16').Value = '3.492.51'
$ws.Range('E16').Value = '  +2.56%  '
$ws.Range('D17').Value = '''7.01'
$ws.Range('E17').Value = '  +6.55%  '
$ws.Range('D18').Value = '2.993.78'
$ws.Range('E18').Value = '  +2.39%  '
$ws.Range('D19').Value = '58.975.89'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = '''427.21'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').Value = '''13.62'
$ws.Range('E21').Value = '  +5.42%  '
$ws.Range('D22').Value = '''0.714'
$ws.Range('E22').Value = '  +6.90%  '
$ws.Range('E23').Value = '  +4.11%  '
$ws.Range('D24').Value = '''13.45'
$ws.Range('E24').Value = '  +4.02%  '
$ws.Range('D25').Value = '''80.50'
$ws.Range('E25').Value = '  +4.15%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('E28').Value = '  +8.00%  '
$ws.Range('D29').Value = '''2.53'
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('D30').Value = '''7.75'
$ws.Range('E30').Value = '  +5.89%  '
$ws.Range('E31').Value = '  +3.67%  '
$ws.Range('D32').Value = '''6.13'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('E33').Value = '  +0.84%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').Value = '''0.972'
$ws.Range('E34').Value = '  +6.35%  '
$ws.Range('B35').Value = 'PEPE'
$ws.Range('C35').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D35').Value = '0.0₃0748'
$ws.Range('E35').Value = '  +18.39%  '
$ws.Range('D36').Value = '''5.76'
$ws.Range('E36').Value = '  +6.56%  '
$ws.Range('E37').Value = '  +4.03%  '
$ws.Range('D38').Value = '''48.88'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('D39').Value = '''8.86'
$ws.Range('E39').Value = '  +4.34%  '
$ws.Range('D40').Value = '''2.74'
$ws.Range('E40').Value = '  +13.20%  '
$ws.Range('D41').Value = '''395.34'
$ws.Range('E41').Value = '  +8.56%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '''0.0349'
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '''0.108'
$ws.Range('E43').Value = '  +0.55%  '
$ws.Range('D44').Value = '2.727.46'
$ws.Range('E44').Value = '  +3.63%  '
$ws.Range('D45').Value = '''0.246'
$ws.Range('E45').Value = '  +6.23%  '
$ws.Range('D46').Value = '''125.53'
$ws.Range('E46').Value = '  +4.97%  '
$ws.Range('E47').Value = '  +0.04%  '
$ws.Range('E48').Value = '  +3.33%  '
$ws.Range('E49').Value = '  +2.31%  '
$ws.Range('D50').Value = '''23.42'
$ws.Range('E50').Value = '  +2.09%  '
$ws.Range('D51').Value = '''32.16'
$ws.Range('E51').Value = '  +16.72%  '
